$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Valor_Observado" column (old column F); shifts G:O left to F:N
$ws.Range("F1").EntireColumn.Delete()

# Rename header cells
$ws.Range("B1").Value = "Config"
$ws.Range("D1").Value = "Dist"
$ws.Range("E1").Value = "Var"

# Reorder / relabel the model-result headers (F1:N1)
$ws.Range("F1").Value = "Block Bootstrapping"
$ws.Range("G1").Value = "Sieve Bootstrap"
$ws.Range("H1").Value = "LSPM"
$ws.Range("I1").Value = "LSPMW"
$ws.Range("J1").Value = "AREPD"
$ws.Range("K1").Value = "MCPS"
$ws.Range("L1").Value = "AV-MCPS"
$ws.Range("M1").Value = "DeepAR"
$ws.Range("N1").Value = "EnCQR-LSTM"

# Update the numeric result data for each scenario row (new simulation run values)
# Row 2
$ws.Range("F2").Value = 0.6369513500959092
$ws.Range("G2").Value = 0.5843994421329943
$ws.Range("H2").Value = 0.570469362795984
$ws.Range("I2").Value = 0.570840921275085
$ws.Range("J2").Value = 0.6533329303549965
$ws.Range("K2").Value = 0.8549203449787076
$ws.Range("L2").Value = 0.7161690003526351
$ws.Range("M2").Value = 0.6283347284026712
$ws.Range("N2").Value = 0.8774877686127051

# Row 3
$ws.Range("F3").Value = 0.7682142964308963
$ws.Range("G3").Value = 0.5875887580573944
$ws.Range("H3").Value = 0.9091143611602889
$ws.Range("I3").Value = 1.093086879111845
$ws.Range("J3").Value = 0.7846454750213733
$ws.Range("K3").Value = 0.6182420795381399
$ws.Range("L3").Value = 0.8092470716527327
$ws.Range("M3").Value = 0.730450496981359
$ws.Range("N3").Value = 0.8793067227491869

# Row 4
$ws.Range("F4").Value = 0.5655028703433915
$ws.Range("G4").Value = 0.5656516677990023
$ws.Range("H4").Value = 0.6579954765286115
$ws.Range("I4").Value = 0.652997466551634
$ws.Range("J4").Value = 0.5706582279067366
$ws.Range("K4").Value = 0.6582565985149048
$ws.Range("L4").Value = 0.6241631250665416
$ws.Range("M4").Value = 0.5673764652125339
$ws.Range("N4").Value = 0.8285478870291892

# Row 5
$ws.Range("F5").Value = 0.5871901013957525
$ws.Range("G5").Value = 0.5654217395649502
$ws.Range("H5").Value = 0.5718470472866304
$ws.Range("I5").Value = 0.5837226012779879
$ws.Range("J5").Value = 0.5875137066800068
$ws.Range("K5").Value = 0.7210538059920016
$ws.Range("L5").Value = 0.5974819202697724
$ws.Range("M5").Value = 0.5979892408691407
$ws.Range("N5").Value = 0.850878739895186

# Row 6
$ws.Range("F6").Value = 0.6491540140440122
$ws.Range("G6").Value = 0.6697446540368869
$ws.Range("H6").Value = 0.6239498082976067
$ws.Range("I6").Value = 0.5680768436293515
$ws.Range("J6").Value = 0.6861917308243309
$ws.Range("K6").Value = 0.6425446912940814
$ws.Range("L6").Value = 0.7035799170340037
$ws.Range("M6").Value = 0.6405699624203186
$ws.Range("N6").Value = 0.8936736347922928

# Row 7
$ws.Range("F7").Value = 0.5894449602282741
$ws.Range("G7").Value = 0.8230185070006231
$ws.Range("H7").Value = 0.586411766223163
$ws.Range("I7").Value = 0.7018933605115637
$ws.Range("J7").Value = 0.5903694533282057
$ws.Range("K7").Value = 0.6938511040699863
$ws.Range("L7").Value = 1.004307895169496
$ws.Range("M7").Value = 0.5871122686892942
$ws.Range("N7").Value = 0.8400430956102833

# Row 8
$ws.Range("F8").Value = 0.866279039966514
$ws.Range("G8").Value = 0.5734746947094902
$ws.Range("H8").Value = 1.389276288467841
$ws.Range("I8").Value = 1.216497922655282
$ws.Range("J8").Value = 0.8693556945192348
$ws.Range("K8").Value = 0.7326250492836474
$ws.Range("L8").Value = 0.9160571857826871
$ws.Range("M8").Value = 0.8188226430512602
$ws.Range("N8").Value = 0.9052844873369854

# Row 9
$ws.Range("F9").Value = 0.7084516708508127
$ws.Range("G9").Value = 0.5852166589566348
$ws.Range("H9").Value = 0.5567462903302045
$ws.Range("I9").Value = 1.030104100435964
$ws.Range("J9").Value = 0.7335859563576617
$ws.Range("K9").Value = 0.7304313654266728
$ws.Range("L9").Value = 0.7194836693236835
$ws.Range("M9").Value = 0.6823386964956435
$ws.Range("N9").Value = 0.8565212200005967

# Row 10
$ws.Range("F10").Value = 0.5919218337762049
$ws.Range("G10").Value = 0.5906312704983701
$ws.Range("H10").Value = 0.6732591212121318
$ws.Range("I10").Value = 0.629151650987046
$ws.Range("J10").Value = 0.5956020494428479
$ws.Range("K10").Value = 0.8982715979309793
$ws.Range("L10").Value = 0.6249703065285882
$ws.Range("M10").Value = 0.5942904237517793
$ws.Range("N10").Value = 0.8440793537434786

# Row 11
$ws.Range("F11").Value = 0.7125432013468992
$ws.Range("G11").Value = 0.5765485154551552
$ws.Range("H11").Value = 0.8241314124589615
$ws.Range("I11").Value = 0.9439419395309855
$ws.Range("J11").Value = 0.6868509931530272
$ws.Range("K11").Value = 0.7833991882374014
$ws.Range("L11").Value = 0.7115120218891509
$ws.Range("M11").Value = 0.6467199027169138
$ws.Range("N11").Value = 0.8431892721908881

# Row 12
$ws.Range("F12").Value = 0.6271611672714793
$ws.Range("G12").Value = 0.5788479493160956
$ws.Range("H12").Value = 0.8814978351959203
$ws.Range("I12").Value = 0.5767767344497744
$ws.Range("J12").Value = 0.6693722492930013
$ws.Range("K12").Value = 0.8133305091970324
$ws.Range("L12").Value = 0.7891414364067056
$ws.Range("M12").Value = 0.6596160004108806
$ws.Range("N12").Value = 0.8859857166476947

# Row 13
$ws.Range("F13").Value = 0.5990583492553717
$ws.Range("G13").Value = 0.5528112836089469
$ws.Range("H13").Value = 0.5600460369563607
$ws.Range("I13").Value = 0.5516873544602805
$ws.Range("J13").Value = 0.6265093932357716
$ws.Range("K13").Value = 1.469510826562523
$ws.Range("L13").Value = 1.288047325402368
$ws.Range("M13").Value = 0.6033556026092327
$ws.Range("N13").Value = 0.873602431543

# Row 14
$ws.Range("F14").Value = 0.5833981083963058
$ws.Range("G14").Value = 0.5762481950672734
$ws.Range("H14").Value = 0.5863248494721115
$ws.Range("I14").Value = 0.5757911765018415
$ws.Range("J14").Value = 0.6074554900455291
$ws.Range("K14").Value = 0.826983955591566
$ws.Range("L14").Value = 0.6633086149499661
$ws.Range("M14").Value = 0.5818056121320607
$ws.Range("N14").Value = 0.8176024116982845

# Row 15
$ws.Range("F15").Value = 0.6332534597162516
$ws.Range("G15").Value = 0.5998194372902358
$ws.Range("H15").Value = 0.6455339243162069
$ws.Range("I15").Value = 0.6554499564311369
$ws.Range("J15").Value = 0.7355551910426364
$ws.Range("K15").Value = 0.8649617968662444
$ws.Range("L15").Value = 0.7112451215453883
$ws.Range("M15").Value = 0.6148501938233483
$ws.Range("N15").Value = 0.8236548418040509

# Row 16
$ws.Range("F16").Value = 0.5670396182421353
$ws.Range("G16").Value = 0.5592416729828991
$ws.Range("H16").Value = 0.577067856635735
$ws.Range("I16").Value = 0.5597833288966785
$ws.Range("J16").Value = 0.5985497539529903
$ws.Range("K16").Value = 0.5881596631715506
$ws.Range("L16").Value = 0.5953917086742714
$ws.Range("M16").Value = 0.5712401996326767
$ws.Range("N16").Value = 0.8005371384073549

# Row 17
$ws.Range("F17").Value = 0.5810647307474653
$ws.Range("G17").Value = 0.5782772633422316
$ws.Range("H17").Value = 0.5797053311861768
$ws.Range("I17").Value = 0.583009934301929
$ws.Range("J17").Value = 0.6358928536473442
$ws.Range("K17").Value = 0.591360892806673
$ws.Range("L17").Value = 0.6261033257984842
$ws.Range("M17").Value = 0.5719116145016903
$ws.Range("N17").Value = 0.7992587313318554

# Row 18
$ws.Range("F18").Value = 0.6362982091857153
$ws.Range("G18").Value = 0.5937815796285062
$ws.Range("H18").Value = 0.6374169392014292
$ws.Range("I18").Value = 0.6199883168474976
$ws.Range("J18").Value = 0.6518736698166593
$ws.Range("K18").Value = 0.5873766550716394
$ws.Range("L18").Value = 0.6170725665584681
$ws.Range("M18").Value = 0.633822776132321
$ws.Range("N18").Value = 0.8422107146169114

# Row 19
$ws.Range("F19").Value = 0.6806667755294935
$ws.Range("G19").Value = 0.7627773915002621
$ws.Range("H19").Value = 0.7392070170497353
$ws.Range("I19").Value = 0.6965323769896112
$ws.Range("J19").Value = 0.7893830486487188
$ws.Range("K19").Value = 0.6841546852991728
$ws.Range("L19").Value = 0.8882132694836258
$ws.Range("M19").Value = 0.670524649442804
$ws.Range("N19").Value = 0.8307266751065397

# Row 20
$ws.Range("F20").Value = 0.7341048996987063
$ws.Range("G20").Value = 0.6418208799089163
$ws.Range("H20").Value = 0.8518809067397505
$ws.Range("I20").Value = 0.7709354435592145
$ws.Range("J20").Value = 0.8853537884206309
$ws.Range("K20").Value = 0.7027551762605483
$ws.Range("L20").Value = 0.7938226937166013
$ws.Range("M20").Value = 0.70355144149167
$ws.Range("N20").Value = 0.8551609176166074

# Row 21
$ws.Range("F21").Value = 0.6144618063131096
$ws.Range("G21").Value = 0.6086739195700983
$ws.Range("H21").Value = 0.7372064946676444
$ws.Range("I21").Value = 0.6096867904141217
$ws.Range("J21").Value = 0.6440304757997999
$ws.Range("K21").Value = 0.6551931917204481
$ws.Range("L21").Value = 0.6984921468127152
$ws.Range("M21").Value = 0.6416968656370493
$ws.Range("N21").Value = 0.8375956833985226

# Row 22
$ws.Range("F22").Value = 0.5649426802959302
$ws.Range("G22").Value = 0.5611258797105646
$ws.Range("H22").Value = 0.5620057222528267
$ws.Range("I22").Value = 0.5675120687147992
$ws.Range("J22").Value = 0.6137617611741494
$ws.Range("K22").Value = 0.7304820535140963
$ws.Range("L22").Value = 0.6191284375737034
$ws.Range("M22").Value = 0.5644383792241516
$ws.Range("N22").Value = 0.7974167652382866

# Row 23
$ws.Range("F23").Value = 0.7159236837559406
$ws.Range("G23").Value = 0.6606679631802898
$ws.Range("H23").Value = 0.6829178424798934
$ws.Range("I23").Value = 0.7432119238151129
$ws.Range("J23").Value = 0.8456817478080344
$ws.Range("K23").Value = 0.7150136676326554
$ws.Range("L23").Value = 0.8453911053076802
$ws.Range("M23").Value = 0.6777267054156471
$ws.Range("N23").Value = 0.8413271756719426

# Row 24
$ws.Range("F24").Value = 0.6322637076424887
$ws.Range("G24").Value = 0.6145203577433802
$ws.Range("H24").Value = 0.6536669564603145
$ws.Range("I24").Value = 0.6204729134081688
$ws.Range("J24").Value = 0.6581745084808737
$ws.Range("K24").Value = 0.9068167041316114
$ws.Range("L24").Value = 0.7748665695067382
$ws.Range("M24").Value = 0.6352118876799946
$ws.Range("N24").Value = 0.8341246528954248

# Row 25
$ws.Range("F25").Value = 0.5889810475808932
$ws.Range("G25").Value = 0.5867506107788001
$ws.Range("H25").Value = 0.5762842947476078
$ws.Range("I25").Value = 0.5803818083840081
$ws.Range("J25").Value = 0.6135909504333112
$ws.Range("K25").Value = 0.6199004693791328
$ws.Range("L25").Value = 0.6003965938727117
$ws.Range("M25").Value = 0.5858136934860649
$ws.Range("N25").Value = 0.8199943522940039

